# Fix Training Data Issue (#48)
# Data was taken from 1 day off due to way NBA stats were shown.
# The "Date" column (BF) held values formatted like "6-22-2012-13";
# correct them to the actual ISO date string "2013-06-22".
#
# NumberFormat is forced to text ("@") before the assignment so Excel's
# COM layer does not silently reinterpret the date-looking string as a
# date serial number; the style is then reset back to Normal so the
# cells keep their original (default) formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("BF2:BF31")
$rng.NumberFormat = "@"
$rng.Value = "2013-06-22"
$rng.Style = "Normal"
